$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.317.08'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.878.96'
$ws.Range("E3").Value = '  -1.74%  '
$ws.Range("E4").Value = '  -0.64%  '
$ws.Range("D5").Value = '245.93'
$ws.Range("E5").Value = '  -4.06%  '
$ws.Range("D6").Value = '0.685'
$ws.Range("E6").Value = '  -6.43%  '
$ws.Range("E7").Value = '  -0.69%  '
$ws.Range("D8").Value = '43.16'
$ws.Range("E8").Value = '  +4.50%  '
$ws.Range("D9").Value = '0.349'
$ws.Range("E9").Value = '  -5.51%  '
$ws.Range("D10").Value = '0.0735'
$ws.Range("E10").Value = '  -3.36%  '
$ws.Range("D11").Value = '0.0967'
$ws.Range("E11").Value = '  -2.29%  '
$ws.Range("D12").Value = '13.03'
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("D13").Value = '2.153.48'
$ws.Range("E13").Value = '  -1.65%  '
$ws.Range("D14").Value = '0.737'
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("D15").Value = '4.92'
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("D16").Value = '1.881.22'
$ws.Range("E16").Value = '  -1.82%  '
$ws.Range("D17").Value = '35.311.88'
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").Value = '73.32'
$ws.Range("E18").Value = '  -2.61%  '
$ws.Range("D19").Value = '0.0₃0819'
$ws.Range("E19").Value = '  -3.42%  '
$ws.Range("D20").Value = '245.15'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").Value = '12.77'
$ws.Range("E21").Value = '  -2.88%  '
$ws.Range("D22").Value = '4.91'
$ws.Range("E22").Value = '  -4.81%  '
$ws.Range("D24").Value = '2.55'
$ws.Range("E24").Value = '  +4.00%  '
$ws.Range("E25").Value = '  -11.59%  '
$ws.Range("D26").Value = '165.33'
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").Value = '8.43'
$ws.Range("E27").Value = '  -3.91%  '
$ws.Range("D28").Value = '18.24'
$ws.Range("E28").Value = '  -3.31%  '
$ws.Range("E29").Value = '  -4.84%  '
$ws.Range("D30").Value = '4.128.47'
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("E31").Value = '  +4.55%  '
$ws.Range("D32").Value = '4.22'
$ws.Range("E32").Value = '  -3.91%  '
$ws.Range("D33").Value = '0.0577'
$ws.Range("E33").Value = '  -2.53%  '
$ws.Range("E34").Value = '  -2.17%  '
$ws.Range("E35").Value = '  -0.65%  '
$ws.Range("D36").Value = '0.846'
$ws.Range("E36").Value = '  -8.00%  '
$ws.Range("D37").Value = '1.97'
$ws.Range("E37").Value = '  -3.74%  '
$ws.Range("D38").Value = '1.58'
$ws.Range("E38").Value = '  -20.84%  '
$ws.Range("E39").Value = '  +7.56%  '
$ws.Range("D40").Value = '96.96'
$ws.Range("E40").Value = '  -0.75%  '
$ws.Range("D41").Value = '16.92'
$ws.Range("E41").Value = '  -1.15%  '
$ws.Range("E42").Value = '  -3.14%  '
$ws.Range("E43").Value = '  -4.53%  '
$ws.Range("D44").Value = '1.286.31'
$ws.Range("E44").Value = '  -4.32%  '
$ws.Range("E45").Value = '  -6.11%  '
$ws.Range("D46").Value = '0.0808'
$ws.Range("E46").Value = '  +6.97%  '
$ws.Range("E48").Value = '  -1.17%  '
$ws.Range("D49").Value = '12.08'
$ws.Range("E49").Value = '  +3.15%  '
$ws.Range("D50").Value = '43.07'
$ws.Range("E50").Value = '  -4.62%  '
$ws.Range("D51").Value = '6.25'
$ws.Range("E51").Value = '  -7.33%  '
